$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3290019035339355
$arr2 = New-Object "object[,]" 1,21
$arr2[0,0] = 115.5762661441167
$arr2[0,1] = 0.003745993474737887
$arr2[0,2] = 0.003311221283276492
$arr2[0,3] = 0.003271002261805447
$arr2[0,4] = 0.002902045145337641
$arr2[0,5] = 0.002856516092974327
$arr2[0,6] = 0.002759359772799551
$arr2[0,7] = 0.002759359772799551
$arr2[0,8] = 0.002737907235398586
$arr2[0,9] = 0.002576774990157842
$arr2[0,10] = 0.002545155628568645
$arr2[0,11] = 0.002533374939058052
$arr2[0,12] = 0.002491612511506542
$arr2[0,13] = 0.002464560822715964
$arr2[0,14] = 0.002326521707540411
$arr2[0,15] = 0.002320505244877328
$arr2[0,16] = 0.002271102864559079
$arr2[0,17] = 0.002271102864559079
$arr2[0,18] = 0.002271102864559079
$arr2[0,19] = 0.002263966386255195
$arr2[0,20] = 0.002252948657780053
$ws.Range("E2:Y2").Value = $arr2

$ws.Range("C3").Value = 0.3210048675537109
$arr3 = New-Object "object[,]" 1,21
$arr3[0,0] = 116.0820356708336
$arr3[0,1] = 0.00390114034510335
$arr3[0,2] = 0.003154378546753798
$arr3[0,3] = 0.003154378546753798
$arr3[0,4] = 0.002984356511183791
$arr3[0,5] = 0.002614503151397367
$arr3[0,6] = 0.002614503151397367
$arr3[0,7] = 0.002614503151397367
$arr3[0,8] = 0.002496781556896464
$arr3[0,9] = 0.002496781556896464
$arr3[0,10] = 0.002496781556896464
$arr3[0,11] = 0.002484546738162722
$arr3[0,12] = 0.002484546738162722
$arr3[0,13] = 0.002484338710470062
$arr3[0,14] = 0.002380461949905595
$arr3[0,15] = 0.002337992552872809
$arr3[0,16] = 0.002337992552872809
$arr3[0,17] = 0.002333518475021782
$arr3[0,18] = 0.00229463697748413
$arr3[0,19] = 0.002267394726359182
$arr3[0,20] = 0.002262807712881747
$ws.Range("E3:Y3").Value = $arr3

$ws.Range("C4").Value = 0.487013578414917
$arr4 = New-Object "object[,]" 1,21
$arr4[0,0] = 115.3677209115485
$arr4[0,1] = 0.003875298282895648
$arr4[0,2] = 0.003188322757866584
$arr4[0,3] = 0.003087980234697744
$arr4[0,4] = 0.003051043328891713
$arr4[0,5] = 0.002951784720873412
$arr4[0,6] = 0.002800067778144085
$arr4[0,7] = 0.002721152937764782
$arr4[0,8] = 0.002504028349357129
$arr4[0,9] = 0.002504028349357129
$arr4[0,10] = 0.002504028349357129
$arr4[0,11] = 0.002504028349357129
$arr4[0,12] = 0.002363532524819749
$arr4[0,13] = 0.002363532524819749
$arr4[0,14] = 0.002363532524819749
$arr4[0,15] = 0.002363532524819749
$arr4[0,16] = 0.002363532524819749
$arr4[0,17] = 0.002277540747599822
$arr4[0,18] = 0.002277540747599822
$arr4[0,19] = 0.002277540747599822
$arr4[0,20] = 0.002248883448568197
$ws.Range("E4:Y4").Value = $arr4

$ws.Range("C5").Value = 0.5479989051818848
$arr5 = New-Object "object[,]" 1,21
$arr5[0,0] = 129.5008367562514
$arr5[0,1] = 0.00390114034510335
$arr5[0,2] = 0.003351749235726567
$arr5[0,3] = 0.003094400686186169
$arr5[0,4] = 0.002971856705439614
$arr5[0,5] = 0.002971856705439614
$arr5[0,6] = 0.002971856705439614
$arr5[0,7] = 0.002909294371298354
$arr5[0,8] = 0.002908298953286547
$arr5[0,9] = 0.002671325835986405
$arr5[0,10] = 0.002627250102949558
$arr5[0,11] = 0.002627250102949558
$arr5[0,12] = 0.002627250102949558
$arr5[0,13] = 0.002627250102949558
$arr5[0,14] = 0.002627250102949558
$arr5[0,15] = 0.002627250102949558
$arr5[0,16] = 0.002596645324756395
$arr5[0,17] = 0.002587793423552649
$arr5[0,18] = 0.002557117657885866
$arr5[0,19] = 0.002546989163288459
$arr5[0,20] = 0.002524382782772931
$ws.Range("E5:Y5").Value = $arr5

$ws.Range("C6").Value = 0.5440025329589844
$arr6 = New-Object "object[,]" 1,21
$arr6[0,0] = 119.1773031544708
$arr6[0,1] = 0.003816966319851988
$arr6[0,2] = 0.003348853838115831
$arr6[0,3] = 0.002965452240544243
$arr6[0,4] = 0.002896068243171647
$arr6[0,5] = 0.002896068243171647
$arr6[0,6] = 0.002762568115126766
$arr6[0,7] = 0.002762568115126766
$arr6[0,8] = 0.002655454454591582
$arr6[0,9] = 0.002578546452246815
$arr6[0,10] = 0.002473197925874135
$arr6[0,11] = 0.002473197925874135
$arr6[0,12] = 0.002473197925874135
$arr6[0,13] = 0.002473197925874135
$arr6[0,14] = 0.002473197925874135
$arr6[0,15] = 0.002450470933164224
$arr6[0,16] = 0.002367915880235332
$arr6[0,17] = 0.002367915880235332
$arr6[0,18] = 0.002367679258896727
$arr6[0,19] = 0.002346915689997923
$arr6[0,20] = 0.002323144311003329
$ws.Range("E6:Y6").Value = $arr6

$ws.Range("C7").Value = 0.4249985218048096
$arr7 = New-Object "object[,]" 1,21
$arr7[0,0] = 117.4314627904896
$arr7[0,1] = 0.00390114034510335
$arr7[0,2] = 0.003106077429802523
$arr7[0,3] = 0.003088420740287512
$arr7[0,4] = 0.003088420740287512
$arr7[0,5] = 0.003012524027192262
$arr7[0,6] = 0.00285937749652403
$arr7[0,7] = 0.00285937749652403
$arr7[0,8] = 0.00285937749652403
$arr7[0,9] = 0.002805399533358707
$arr7[0,10] = 0.002678906427032384
$arr7[0,11] = 0.002643496840581979
$arr7[0,12] = 0.002643496840581979
$arr7[0,13] = 0.002587195590327631
$arr7[0,14] = 0.002511841649580642
$arr7[0,15] = 0.002464887513052288
$arr7[0,16] = 0.002350473465424561
$arr7[0,17] = 0.002350473465424561
$arr7[0,18] = 0.002292844280938532
$arr7[0,19] = 0.002292844280938532
$arr7[0,20] = 0.002289112335097263
$ws.Range("E7:Y7").Value = $arr7

$ws.Range("C8").Value = 0.4049961566925049
$arr8 = New-Object "object[,]" 1,21
$arr8[0,0] = 111.821973236918
$arr8[0,1] = 0.003590901701430972
$arr8[0,2] = 0.003347972639388806
$arr8[0,3] = 0.003152854258406206
$arr8[0,4] = 0.003037416513677748
$arr8[0,5] = 0.003037416513677748
$arr8[0,6] = 0.002888145139176945
$arr8[0,7] = 0.002879339649445388
$arr8[0,8] = 0.002738439701221337
$arr8[0,9] = 0.002709927482098853
$arr8[0,10] = 0.002524346503001655
$arr8[0,11] = 0.002524346503001655
$arr8[0,12] = 0.002493418117328882
$arr8[0,13] = 0.002385573940101471
$arr8[0,14] = 0.002321587206784032
$arr8[0,15] = 0.002252201359514283
$arr8[0,16] = 0.002252201359514283
$arr8[0,17] = 0.002246477327907021
$arr8[0,18] = 0.002246477327907021
$arr8[0,19] = 0.002179765560173839
$arr8[0,20] = 0.002179765560173839
$ws.Range("E8:Y8").Value = $arr8

$ws.Range("C9").Value = 0.4349937438964844
$arr9 = New-Object "object[,]" 1,21
$arr9[0,0] = 122.3932759372856
$arr9[0,1] = 0.00390114034510335
$arr9[0,2] = 0.003255092378214194
$arr9[0,3] = 0.003076931133893406
$arr9[0,4] = 0.002946731242252434
$arr9[0,5] = 0.002946731242252434
$arr9[0,6] = 0.002946731242252434
$arr9[0,7] = 0.002946731242252434
$arr9[0,8] = 0.002923845474666755
$arr9[0,9] = 0.002693520465047228
$arr9[0,10] = 0.002693520465047228
$arr9[0,11] = 0.002693520465047228
$arr9[0,12] = 0.002693520465047228
$arr9[0,13] = 0.002693520465047228
$arr9[0,14] = 0.002678015235453794
$arr9[0,15] = 0.002566830535159747
$arr9[0,16] = 0.002562960805302827
$arr9[0,17] = 0.002521946752101967
$arr9[0,18] = 0.002499489200416713
$arr9[0,19] = 0.002429860868424266
$arr9[0,20] = 0.002385833838933442
$ws.Range("E9:Y9").Value = $arr9

$ws.Range("C10").Value = 0.4180018901824951
$arr10 = New-Object "object[,]" 1,21
$arr10[0,0] = 121.394236595961
$arr10[0,1] = 0.00390114034510335
$arr10[0,2] = 0.003562289639906234
$arr10[0,3] = 0.003181439551302454
$arr10[0,4] = 0.003086772112001411
$arr10[0,5] = 0.002888602871172946
$arr10[0,6] = 0.002888602871172946
$arr10[0,7] = 0.002879350245184433
$arr10[0,8] = 0.002746191571010877
$arr10[0,9] = 0.002746191571010877
$arr10[0,10] = 0.002746191571010877
$arr10[0,11] = 0.002687363836138241
$arr10[0,12] = 0.002657060311746683
$arr10[0,13] = 0.002570654203022126
$arr10[0,14] = 0.002570654203022126
$arr10[0,15] = 0.002486518986518995
$arr10[0,16] = 0.002469341334502613
$arr10[0,17] = 0.002469341334502613
$arr10[0,18] = 0.002404124137270383
$arr10[0,19] = 0.002392639939210207
$arr10[0,20] = 0.002366359387835497
$ws.Range("E10:Y10").Value = $arr10

$ws.Range("C11").Value = 0.3959970474243164
$arr11 = New-Object "object[,]" 1,21
$arr11[0,0] = 115.0059951187304
$arr11[0,1] = 0.003674671469616441
$arr11[0,2] = 0.003428591630406239
$arr11[0,3] = 0.003198387528082945
$arr11[0,4] = 0.003198387528082945
$arr11[0,5] = 0.003198387528082945
$arr11[0,6] = 0.002882516527010096
$arr11[0,7] = 0.002882516527010096
$arr11[0,8] = 0.002792470487521237
$arr11[0,9] = 0.002731008768422358
$arr11[0,10] = 0.002731008768422358
$arr11[0,11] = 0.002695263078483754
$arr11[0,12] = 0.002559899364851555
$arr11[0,13] = 0.002559899364851555
$arr11[0,14] = 0.002559899364851555
$arr11[0,15] = 0.002435825442533147
$arr11[0,16] = 0.002405991964749855
$arr11[0,17] = 0.002370610493854384
$arr11[0,18] = 0.002334412140038441
$arr11[0,19] = 0.002305415447388806
$arr11[0,20] = 0.002241832263523009
$ws.Range("E11:Y11").Value = $arr11
Write-Output "done"
